$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    'SELECT "id", "customer_number", "customer_name", "client_id", "client_name", "captain_name", "delivery_date", "order_status", "shop_to_delivery_km", "order_created_at", "order_accepted_at", "start_ride_at", "reached_shop_at", "order_picked_at", "shipped_at", "reached_dest_at", "final_status_at", "cancellation_reason" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "order_created_at" >= (CURRENT_TIMESTAMP - INTERVAL ''1 year'') ORDER BY "order_created_at" DESC LIMIT 10;',
    'SELECT MAX("order_created_at") AS "Last_Order_Time" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'';',
    'SELECT "delivery_date" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "delivery_date" = (SELECT MAX("delivery_date") FROM "updated_table" WHERE "client_name" = ''MC DONALDS'');',
    'SELECT "id", "customer_number", "customer_name", "client_id", "client_name", "captain_name", "delivery_date", "order_status", "shop_to_delivery_km", "order_created_at", "order_accepted_at", "start_ride_at", "reached_shop_at", "order_picked_at", "shipped_at", "reached_dest_at", "final_status_at", "cancellation_reason" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "order_created_at" <= ''2025-09-12 12:37:40.233648'' ORDER BY "order_created_at" DESC LIMIT 10;',
    'SELECT "id", "customer_number", "customer_name", "client_id", "client_name", "captain_name", "delivery_date", "order_status", "shop_to_delivery_km", "order_created_at", "order_accepted_at", "start_ride_at", "reached_shop_at", "order_picked_at", "shipped_at", "reached_dest_at", "final_status_at", "cancellation_reason" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "order_created_at" <= ''2025-09-12 12:39:10.848770'' ORDER BY "order_created_at" DESC LIMIT 10;',
    'SELECT "order_created_at" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "order_created_at" = (SELECT MAX("order_created_at") FROM "updated_table" WHERE "client_name" = ''MC DONALDS'');'
)

$startRow = 8
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i]
}

$wb.Save()
